# Auto-generated Excel COM-interop script
# Applies numeric cell updates across sheets ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR
# to match the target OOXML diff (scheduled Sheets runner data refresh).

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 1308180.8
$ws.Range("J17").Value = 1337140.2
$ws.Range("L17").Value = 4011420.6
$ws.Range("N17").Value = -4011756.6
$ws.Range("H18").Value = 966.9535
$ws.Range("I18").Value = 966.9535
$ws.Range("K18").Value = 966.9535
$ws.Range("M18").Value = -682.9535
$ws.Range("H33").Value = 16672722
$ws.Range("J33").Value = 8916.625
$ws.Range("L33").Value = 8916.625
$ws.Range("N33").Value = -9374.625
$ws.Range("H38").Value = 2742.7778
$ws.Range("I38").Value = 396.66666
$ws.Range("K38").Value = 1189.99998
$ws.Range("M38").Value = -817.9999800000001
$ws.Range("H39").Value = 13193.625
$ws.Range("I39").Value = 758.3333
$ws.Range("K39").Value = 2274.9999
$ws.Range("M39").Value = -1978.9999
$ws.Range("H40").Value = 3999.6428
$ws.Range("I40").Value = 7624.75
$ws.Range("J40").Value = 2549.6
$ws.Range("K40").Value = 7624.75
$ws.Range("L40").Value = 2549.6
$ws.Range("M40").Value = -7449.75
$ws.Range("N40").Value = -2899.6
$ws.Range("H43").Value = 2757.6667
$ws.Range("I43").Value = 2514
$ws.Range("J43").Value = 3732.3333
$ws.Range("K43").Value = 2514
$ws.Range("L43").Value = 3732.3333
$ws.Range("M43").Value = -2445
$ws.Range("N43").Value = -3870.3333
$ws.Range("H62").Value = 9588.4
$ws.Range("I62").Value = 2879.6
$ws.Range("K62").Value = 2879.6
$ws.Range("M62").Value = -2255.6
$ws.Range("H65").Value = 9588.4
$ws.Range("I65").Value = 2879.6
$ws.Range("K65").Value = 14398
$ws.Range("M65").Value = -11278
$ws.Range("H95").Value = 48021.9
$ws.Range("J95").Value = 47252.11
$ws.Range("L95").Value = 47252.11
$ws.Range("N95").Value = -52744.11
$ws.Range("H96").Value = 1467.9375
$ws.Range("I96").Value = 1491.1154
$ws.Range("J96").Value = 1367.5
$ws.Range("K96").Value = 4473.3462
$ws.Range("L96").Value = 4102.5
$ws.Range("M96").Value = -3100.3462
$ws.Range("N96").Value = -6848.5
$ws.Range("H98").Value = 3396.0977
$ws.Range("I98").Value = 2811.6775
$ws.Range("K98").Value = 2811.6775
$ws.Range("M98").Value = -1313.6775
$ws.Range("H100").Value = 6348.5557
$ws.Range("I100").Value = 4467.4
$ws.Range("J100").Value = 8700
$ws.Range("K100").Value = 4467.4
$ws.Range("L100").Value = 8700
$ws.Range("M100").Value = -3926.4
$ws.Range("N100").Value = -9782
$ws.Range("H101").Value = 1996.4445
$ws.Range("I101").Value = 1828
$ws.Range("K101").Value = 5484
$ws.Range("M101").Value = -3862
$ws.Range("H106").Value = 6548.737
$ws.Range("I106").Value = 4488.7856
$ws.Range("J106").Value = 12316.6
$ws.Range("K106").Value = 4488.7856
$ws.Range("L106").Value = 12316.6
$ws.Range("M106").Value = -3857.7856
$ws.Range("N106").Value = -13578.6
$ws.Range("H111").Value = 2106.111
$ws.Range("I111").Value = 800
$ws.Range("K111").Value = 2400
$ws.Range("M111").Value = 667
$ws.Range("H113").Value = 14866.25
$ws.Range("I113").Value = 14866.25
$ws.Range("K113").Value = 14866.25
$ws.Range("M113").Value = -11612.25
$ws.Range("H116").Value = 5600.227
$ws.Range("I116").Value = 5595.2666
$ws.Range("J116").Value = 5610.857
$ws.Range("K116").Value = 5595.2666
$ws.Range("L116").Value = 5610.857
$ws.Range("M116").Value = -2153.2666
$ws.Range("N116").Value = -12494.857
$ws.Range("H122").Value = 3396.0977
$ws.Range("I122").Value = 2811.6775
$ws.Range("K122").Value = 8435.032499999999
$ws.Range("M122").Value = -5985.032499999999
$ws.Range("H125").Value = 5135.143
$ws.Range("I125").Value = 6662
$ws.Range("J125").Value = 1318
$ws.Range("K125").Value = 59958
$ws.Range("L125").Value = 11862
$ws.Range("M125").Value = -57498
$ws.Range("N125").Value = -16782
$ws.Range("H131").Value = 2217.0833
$ws.Range("I131").Value = 1691.8182
$ws.Range("K131").Value = 5075.4546
$ws.Range("M131").Value = -35.45460000000003
$ws.Range("H133").Value = 34633
$ws.Range("J133").Value = 34633
$ws.Range("L133").Value = 34633
$ws.Range("N133").Value = -44753
$ws.Range("H134").Value = 70824.28999999999
$ws.Range("J134").Value = 70824.28999999999
$ws.Range("L134").Value = 70824.28999999999
$ws.Range("N134").Value = -80964.28999999999
$ws.Range("H136").Value = 87741.25
$ws.Range("J136").Value = 87741.25
$ws.Range("L136").Value = 87741.25
$ws.Range("N136").Value = -97941.25
$ws.Range("H137").Value = 51329.668
$ws.Range("I137").Value = 2656.3333
$ws.Range("J137").Value = 100003
$ws.Range("K137").Value = 7968.999899999999
$ws.Range("L137").Value = 300009
$ws.Range("M137").Value = -5418.999899999999
$ws.Range("N137").Value = -305109
$ws.Range("H138").Value = 2323.08
$ws.Range("I138").Value = 1764.2084
$ws.Range("J138").Value = 2499.5657
$ws.Range("K138").Value = 5292.6252
$ws.Range("L138").Value = 7498.6971
$ws.Range("M138").Value = -152.6252000000004
$ws.Range("N138").Value = -17778.6971
$ws.Range("H139").Value = 0
$ws.Range("J139").Value = 0
$ws.Range("L139").Value = 0
$ws.Range("N139").ClearContents()

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 22999.625
$ws.Range("I2").Value = 17332.834
$ws.Range("J2").Value = 40000
$ws.Range("K2").Value = 17332.834
$ws.Range("L2").Value = 40000
$ws.Range("M2").Value = -17219.834
$ws.Range("N2").Value = -40226
$ws.Range("H32").Value = 4943.7236
$ws.Range("I32").Value = 2867.9836
$ws.Range("J32").Value = 13385.066
$ws.Range("K32").Value = 2867.9836
$ws.Range("L32").Value = 13385.066
$ws.Range("M32").Value = -2580.9836
$ws.Range("N32").Value = -13959.066
$ws.Range("H33").Value = 9500
$ws.Range("I33").Value = 9500
$ws.Range("K33").Value = 9500
$ws.Range("M33").Value = -9171
$ws.Range("H45").Value = 6181.2856
$ws.Range("I45").Value = 5730.2666
$ws.Range("J45").Value = 7308.8335
$ws.Range("K45").Value = 5730.2666
$ws.Range("L45").Value = 7308.8335
$ws.Range("M45").Value = -5353.2666
$ws.Range("N45").Value = -8062.8335
$ws.Range("H61").Value = 52314.723
$ws.Range("I61").Value = 2539.7812
$ws.Range("K61").Value = 2539.7812
$ws.Range("M61").Value = -2327.7812
$ws.Range("H74").Value = 12091.053
$ws.Range("I74").Value = 1728.7587
$ws.Range("K74").Value = 1728.7587
$ws.Range("M74").Value = -854.7587000000001
$ws.Range("H77").Value = 12091.053
$ws.Range("I77").Value = 1728.7587
$ws.Range("K77").Value = 8643.7935
$ws.Range("M77").Value = -4275.7935
$ws.Range("H95").Value = 7538.6665
$ws.Range("J95").Value = 7538.6665
$ws.Range("L95").Value = 7538.6665
$ws.Range("N95").Value = -13030.6665
$ws.Range("H96").Value = 0
$ws.Range("J96").Value = 0
$ws.Range("L96").Value = 0
$ws.Range("N96").ClearContents()
$ws.Range("H105").Value = 101496
$ws.Range("J105").Value = 101496
$ws.Range("L105").Value = 101496
$ws.Range("N105").Value = -108484
$ws.Range("H106").Value = 180000
$ws.Range("J106").Value = 180000
$ws.Range("L106").Value = 180000
$ws.Range("N106").Value = -182524
$ws.Range("H110").Value = 5689856
$ws.Range("I110").Value = 6998670
$ws.Range("K110").Value = 6998670
$ws.Range("M110").Value = -6996625
$ws.Range("H116").Value = 22999.625
$ws.Range("I116").Value = 17332.834
$ws.Range("J116").Value = 40000
$ws.Range("K116").Value = 17332.834
$ws.Range("L116").Value = 40000
$ws.Range("M116").Value = -15038.834
$ws.Range("N116").Value = -44588
$ws.Range("H120").Value = 69900
$ws.Range("J120").Value = 69900
$ws.Range("L120").Value = 69900
$ws.Range("N120").Value = -79576
$ws.Range("H122").Value = 551491
$ws.Range("I122").Value = 862070.6
$ws.Range("K122").Value = 2586211.8
$ws.Range("M122").Value = -2583761.8
$ws.Range("H132").Value = 13343.842
$ws.Range("I132").Value = 1770.8846
$ws.Range("J132").Value = 38418.582
$ws.Range("K132").Value = 5312.6538
$ws.Range("L132").Value = 115255.746
$ws.Range("M132").Value = -2782.6538
$ws.Range("N132").Value = -120315.746
$ws.Range("H136").Value = 52314.723
$ws.Range("I136").Value = 2539.7812
$ws.Range("K136").Value = 7619.3436
$ws.Range("M136").Value = -5069.3436
$ws.Range("H137").Value = 46544.91
$ws.Range("I137").Value = 40000
$ws.Range("J137").Value = 48999.25
$ws.Range("K137").Value = 40000
$ws.Range("L137").Value = 48999.25
$ws.Range("M137").Value = -34900
$ws.Range("N137").Value = -59199.25

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 22999.625
$ws.Range("I3").Value = 17332.834
$ws.Range("J3").Value = 40000
$ws.Range("K3").Value = 17332.834
$ws.Range("L3").Value = 40000
$ws.Range("M3").Value = -17218.834
$ws.Range("N3").Value = -40228
$ws.Range("H20").Value = 33992.04
$ws.Range("I20").Value = 10799.363
$ws.Range("K20").Value = 10799.363
$ws.Range("M20").Value = -10552.363
$ws.Range("H62").Value = 89181
$ws.Range("J62").Value = 89181
$ws.Range("L62").Value = 89181
$ws.Range("N62").Value = -90553
$ws.Range("H65").Value = 89181
$ws.Range("J65").Value = 89181
$ws.Range("L65").Value = 267543
$ws.Range("N65").Value = -274407
$ws.Range("H76").Value = 1557
$ws.Range("J76").Value = 1557
$ws.Range("L76").Value = 1557
$ws.Range("N76").Value = -2187
$ws.Range("H79").Value = 1557
$ws.Range("J79").Value = 1557
$ws.Range("L79").Value = 1557
$ws.Range("N79").Value = -3741
$ws.Range("H86").Value = 1730.5128
$ws.Range("I86").Value = 1524.3928
$ws.Range("J86").Value = 2255.182
$ws.Range("K86").Value = 1524.3928
$ws.Range("L86").Value = 2255.182
$ws.Range("M86").Value = -401.3928000000001
$ws.Range("N86").Value = -4501.182
$ws.Range("H89").Value = 1730.5128
$ws.Range("I89").Value = 1524.3928
$ws.Range("J89").Value = 2255.182
$ws.Range("K89").Value = 7621.964
$ws.Range("L89").Value = 11275.91
$ws.Range("M89").Value = -2005.964
$ws.Range("N89").Value = -22507.91
$ws.Range("H99").Value = 30816.889
$ws.Range("I99").Value = 31436.715
$ws.Range("J99").Value = 28647.5
$ws.Range("K99").Value = 31436.715
$ws.Range("L99").Value = 28647.5
$ws.Range("M99").Value = -29938.715
$ws.Range("N99").Value = -31643.5
$ws.Range("H103").Value = 5850
$ws.Range("J103").Value = 5850
$ws.Range("L103").Value = 5850
$ws.Range("N103").Value = -8194
$ws.Range("H105").Value = 2457
$ws.Range("I105").Value = 2358.2727
$ws.Range("K105").Value = 2358.2727
$ws.Range("M105").Value = -611.2727
$ws.Range("H134").Value = 53693.88
$ws.Range("I134").Value = 68488.81
$ws.Range("J134").Value = 27391.777
$ws.Range("K134").Value = 205466.43
$ws.Range("L134").Value = 82175.33099999999
$ws.Range("M134").Value = -202931.43
$ws.Range("N134").Value = -87245.33099999999

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H8").Value = 13505
$ws.Range("I8").Value = 7000
$ws.Range("K8").Value = 7000
$ws.Range("M8").Value = -6860
$ws.Range("H16").Value = 10241.9
$ws.Range("I16").Value = 6101.5835
$ws.Range("K16").Value = 6101.5835
$ws.Range("M16").Value = -5814.5835
$ws.Range("H22").Value = 835.5769
$ws.Range("I22").Value = 412.10526
$ws.Range("J22").Value = 1985
$ws.Range("K22").Value = 412.10526
$ws.Range("L22").Value = 1985
$ws.Range("M22").Value = -62.10525999999999
$ws.Range("N22").Value = -2685
$ws.Range("H31").Value = 22557.572
$ws.Range("I31").Value = 19634.334
$ws.Range("J31").Value = 23726.867
$ws.Range("K31").Value = 19634.334
$ws.Range("L31").Value = 23726.867
$ws.Range("M31").Value = -19339.334
$ws.Range("N31").Value = -24316.867
$ws.Range("H33").Value = 3368.5
$ws.Range("I33").Value = 1158
$ws.Range("J33").Value = 10000
$ws.Range("K33").Value = 1158
$ws.Range("L33").Value = 10000
$ws.Range("M33").Value = -779
$ws.Range("N33").Value = -10758
$ws.Range("H34").Value = 22557.572
$ws.Range("I34").Value = 19634.334
$ws.Range("J34").Value = 23726.867
$ws.Range("K34").Value = 19634.334
$ws.Range("L34").Value = 23726.867
$ws.Range("M34").Value = -19432.334
$ws.Range("N34").Value = -24130.867
$ws.Range("H58").Value = 15696.436
$ws.Range("I58").Value = 4348.091
$ws.Range("K58").Value = 4348.091
$ws.Range("M58").Value = -4145.091
$ws.Range("H94").Value = 4988.643
$ws.Range("J94").Value = 812.5
$ws.Range("L94").Value = 812.5
$ws.Range("N94").Value = -1714.5
$ws.Range("H105").Value = 13813.5
$ws.Range("I105").Value = 17727.334
$ws.Range("J105").Value = 9899.666999999999
$ws.Range("K105").Value = 17727.334
$ws.Range("L105").Value = 9899.666999999999
$ws.Range("M105").Value = -15980.334
$ws.Range("N105").Value = -13393.667
$ws.Range("H107").Value = 2062.2
$ws.Range("I107").Value = 511.85715
$ws.Range("J107").Value = 3095.762
$ws.Range("K107").Value = 511.85715
$ws.Range("L107").Value = 3095.762
$ws.Range("M107").Value = 1408.14285
$ws.Range("N107").Value = -6935.762000000001
$ws.Range("H113").Value = 10241.9
$ws.Range("I113").Value = 6101.5835
$ws.Range("K113").Value = 6101.5835
$ws.Range("M113").Value = -3931.5835
$ws.Range("H122").Value = 7201.5
$ws.Range("I122").Value = 1049.75
$ws.Range("K122").Value = 3149.25
$ws.Range("M122").Value = -699.25
$ws.Range("H125").Value = 84663
$ws.Range("J125").Value = 84663
$ws.Range("L125").Value = 84663
$ws.Range("N125").Value = -89583
$ws.Range("H132").Value = 6185.033
$ws.Range("J132").Value = 12557.583
$ws.Range("L132").Value = 37672.749
$ws.Range("N132").Value = -42732.749
$ws.Range("H134").Value = 4392.2
$ws.Range("I134").Value = 1845.5676
$ws.Range("J134").Value = 11640.308
$ws.Range("K134").Value = 5536.7028
$ws.Range("L134").Value = 34920.924
$ws.Range("M134").Value = -3001.7028
$ws.Range("N134").Value = -39990.924
$ws.Range("H136").Value = 15696.436
$ws.Range("I136").Value = 4348.091
$ws.Range("K136").Value = 13044.273
$ws.Range("M136").Value = -10494.273

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 3228622.8
$ws.Range("I4").Value = 3400055.2
$ws.Range("K4").Value = 10200165.6
$ws.Range("M4").Value = -10200053.6
$ws.Range("H5").Value = 8929855
$ws.Range("I5").Value = 1296
$ws.Range("K5").Value = 3888
$ws.Range("M5").Value = -3776
$ws.Range("H8").Value = 10077152
$ws.Range("I8").Value = 10077152
$ws.Range("K8").Value = 30231456
$ws.Range("M8").Value = -30231317
$ws.Range("H92").Value = 866
$ws.Range("J92").Value = 949
$ws.Range("L92").Value = 2847
$ws.Range("N92").Value = -5343
$ws.Range("H107").Value = 860.46155
$ws.Range("I107").Value = 561
$ws.Range("J107").Value = 1159.9231
$ws.Range("K107").Value = 1683
$ws.Range("L107").Value = 3479.7693
$ws.Range("M107").Value = 237
$ws.Range("N107").Value = -7319.7693
$ws.Range("H122").Value = 11959780
$ws.Range("J122").Value = 2366792.5
$ws.Range("L122").Value = 21301132.5
$ws.Range("N122").Value = -21306032.5
$ws.Range("H131").Value = 1407.0605
$ws.Range("I131").Value = 1024.3334
$ws.Range("J131").Value = 1445.3334
$ws.Range("K131").Value = 3073.0002
$ws.Range("L131").Value = 4336.0002
$ws.Range("M131").Value = 1966.9998
$ws.Range("N131").Value = -14416.0002
$ws.Range("H135").Value = 8929855
$ws.Range("I135").Value = 1296
$ws.Range("K135").Value = 11664
$ws.Range("M135").Value = -9129

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H29").Value = 5997.5
$ws.Range("I29").Value = 2000
$ws.Range("J29").Value = 9995
$ws.Range("K29").Value = 2000
$ws.Range("L29").Value = 9995
$ws.Range("M29").Value = -1710
$ws.Range("N29").Value = -10575
$ws.Range("H70").Value = 19410.045
$ws.Range("I70").Value = 16954.416
$ws.Range("J70").Value = 22356.8
$ws.Range("K70").Value = 16954.416
$ws.Range("L70").Value = 22356.8
$ws.Range("M70").Value = -16684.416
$ws.Range("N70").Value = -22896.8
$ws.Range("H73").Value = 19410.045
$ws.Range("I73").Value = 16954.416
$ws.Range("J73").Value = 22356.8
$ws.Range("K73").Value = 16954.416
$ws.Range("L73").Value = 22356.8
$ws.Range("M73").Value = -16018.416
$ws.Range("N73").Value = -24228.8
$ws.Range("H88").Value = 180497.5
$ws.Range("J88").Value = 180497.5
$ws.Range("L88").Value = 180497.5
$ws.Range("N88").Value = -181399.5
$ws.Range("H91").Value = 180497.5
$ws.Range("J91").Value = 180497.5
$ws.Range("L91").Value = 180497.5
$ws.Range("N91").Value = -183617.5
$ws.Range("H99").Value = 5473.875
$ws.Range("I99").Value = 3113
$ws.Range("K99").Value = 3113
$ws.Range("M99").Value = -867
$ws.Range("H102").Value = 352294.88
$ws.Range("I102").Value = 441667.56
$ws.Range("K102").Value = 441667.56
$ws.Range("M102").Value = -440045.56
$ws.Range("H126").Value = 831106.3
$ws.Range("I126").Value = 2196475.2
$ws.Range("K126").Value = 6589425.600000001
$ws.Range("M126").Value = -6586955.600000001
$ws.Range("H132").Value = 10014.628
$ws.Range("I132").Value = 8749.343000000001
$ws.Range("J132").Value = 15550.25
$ws.Range("K132").Value = 26248.029
$ws.Range("L132").Value = 46650.75
$ws.Range("M132").Value = -23718.029
$ws.Range("N132").Value = -51710.75
$ws.Range("H139").Value = 85000
$ws.Range("J139").Value = 85000
$ws.Range("L139").Value = 85000
$ws.Range("N139").Value = -95280

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 100003290
$ws.Range("I16").Value = 125003610
$ws.Range("K16").Value = 125003610
$ws.Range("M16").Value = -125003440
$ws.Range("H40").Value = 1157978.9
$ws.Range("I40").Value = 2824.975
$ws.Range("J40").Value = 5358538.5
$ws.Range("K40").Value = 2824.975
$ws.Range("L40").Value = 5358538.5
$ws.Range("M40").Value = -2688.975
$ws.Range("N40").Value = -5358810.5
$ws.Range("H46").Value = 3136.1667
$ws.Range("I46").Value = 2963.4
$ws.Range("K46").Value = 2963.4
$ws.Range("M46").Value = -2775.4
$ws.Range("H93").Value = 29419344
$ws.Range("I93").Value = 43483410
$ws.Range("J93").Value = 12667.363
$ws.Range("K93").Value = 43483410
$ws.Range("L93").Value = 12667.363
$ws.Range("M93").Value = -43482162
$ws.Range("N93").Value = -15163.363
$ws.Range("H100").Value = 6406.3184
$ws.Range("I100").Value = 4408.9414
$ws.Range("J100").Value = 13197.4
$ws.Range("K100").Value = 4408.9414
$ws.Range("L100").Value = 13197.4
$ws.Range("M100").Value = -3867.9414
$ws.Range("N100").Value = -14279.4
$ws.Range("H106").Value = 15724.75
$ws.Range("J106").Value = 15724.75
$ws.Range("L106").Value = 15724.75
$ws.Range("N106").Value = -18248.75
$ws.Range("H122").Value = 45464160
$ws.Range("I122").Value = 71435840
$ws.Range("J122").Value = 13712.25
$ws.Range("K122").Value = 214307520
$ws.Range("L122").Value = 41136.75
$ws.Range("M122").Value = -214305070
$ws.Range("N122").Value = -46036.75
$ws.Range("H132").Value = 8935.948
$ws.Range("I132").Value = 3827.3845
$ws.Range("J132").Value = 19153.076
$ws.Range("K132").Value = 11482.1535
$ws.Range("L132").Value = 57459.228
$ws.Range("M132").Value = -8952.1535
$ws.Range("N132").Value = -62519.228
$ws.Range("H136").Value = 40434.59
$ws.Range("I136").Value = 41706.84
$ws.Range("J136").Value = 29832.5
$ws.Range("K136").Value = 125120.52
$ws.Range("L136").Value = 89497.5
$ws.Range("M136").Value = -122570.52
$ws.Range("N136").Value = -94597.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H32").Value = 10257.25
$ws.Range("J32").Value = 30029
$ws.Range("L32").Value = 30029
$ws.Range("N32").Value = -30663
$ws.Range("H74").Value = 22582.5
$ws.Range("J74").Value = 22582.5
$ws.Range("L74").Value = 22582.5
$ws.Range("N74").Value = -24454.5
$ws.Range("H77").Value = 22582.5
$ws.Range("J77").Value = 22582.5
$ws.Range("L77").Value = 67747.5
$ws.Range("N77").Value = -77107.5
$ws.Range("H97").Value = 46183.75
$ws.Range("J97").Value = 46183.75
$ws.Range("L97").Value = 46183.75
$ws.Range("N97").Value = -48165.75
$ws.Range("H100").Value = 523.2222
$ws.Range("I100").Value = 441.33334
$ws.Range("J100").Value = 932.6667
$ws.Range("K100").Value = 882.66668
$ws.Range("L100").Value = 1865.3334
$ws.Range("M100").Value = -341.66668
$ws.Range("N100").Value = -2947.3334
$ws.Range("H104").Value = 16000
$ws.Range("J104").Value = 16000
$ws.Range("L104").Value = 16000
$ws.Range("N104").Value = -22988
$ws.Range("H107").Value = 2942.1667
$ws.Range("I107").Value = 869.2353000000001
$ws.Range("K107").Value = 2607.7059
$ws.Range("M107").Value = -687.7058999999999
$ws.Range("H113").Value = 1256.25
$ws.Range("I113").Value = 516.25
$ws.Range("J113").Value = 1996.25
$ws.Range("K113").Value = 1548.75
$ws.Range("L113").Value = 5988.75
$ws.Range("M113").Value = 621.25
$ws.Range("N113").Value = -10328.75
$ws.Range("H126").Value = 11119484
$ws.Range("I126").Value = 2451.1667
$ws.Range("J126").Value = 33353550
$ws.Range("K126").Value = 7353.500100000001
$ws.Range("L126").Value = 100060650
$ws.Range("M126").Value = -4883.500100000001
$ws.Range("N126").Value = -100065590
$ws.Range("H129").Value = 22666.666
$ws.Range("J129").Value = 22666.666
$ws.Range("L129").Value = 22666.666
$ws.Range("N129").Value = -32666.666
$ws.Range("H132").Value = 6720.7827
$ws.Range("I132").Value = 1737.5769
$ws.Range("J132").Value = 13198.95
$ws.Range("K132").Value = 5212.7307
$ws.Range("L132").Value = 39596.85000000001
$ws.Range("M132").Value = -2682.7307
$ws.Range("N132").Value = -44656.85000000001
$ws.Range("H136").Value = 6772.3335
$ws.Range("I136").Value = 700.2273
$ws.Range("K136").Value = 2100.6819
$ws.Range("M136").Value = 449.3181
$ws.Range("H141").Value = 73999.336
$ws.Range("J141").Value = 73999.336
$ws.Range("L141").Value = 73999.336
$ws.Range("N141").Value = -84359.336
